$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-converted to a number
# by Excel's type inference; force Text format first so they round-trip as strings.
$textCells = @("D5", "D6", "D10", "D11", "D12", "D13", "D14", "D15", "D18", "D20", "D21", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D37", "D39", "D40", "D41", "D42", "D43", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "62.582.91"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "3.023.17"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "579.79"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").Value = "148.90"
$ws.Range("E6").Value = "  -3.55%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("D9").Value = "3.023.93"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").Value = "0.151"
$ws.Range("E10").Value = "  -4.03%  "
$ws.Range("D11").Value = "5.67"
$ws.Range("D12").Value = "0.443"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  -3.94%  "
$ws.Range("D14").Value = "35.32"
$ws.Range("E14").Value = "  -4.90%  "
$ws.Range("D15").Value = "0.122"
$ws.Range("D16").Value = "3.526.39"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "62.578.47"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "7.01"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").Value = "3.020.22"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").Value = "468.17"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").Value = "14.00"
$ws.Range("E21").Value = "  -3.35%  "
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").Value = "80.95"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "12.48"
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("D27").Value = "10.47"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "7.23"
$ws.Range("E30").Value = "  -2.93%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.62"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "2.16"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "27.30"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("E34").Value = "  -5.52%  "
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "0.0₃0795"
$ws.Range("E36").Value = "  -6.40%  "
$ws.Range("D37").Value = "5.78"
$ws.Range("E37").Value = "  -4.21%  "
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").Value = "50.28"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "9.02"
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "2.97"
$ws.Range("E41").Value = "  -14.57%  "
$ws.Range("D42").Value = "424.26"
$ws.Range("E42").Value = "  -5.83%  "
$ws.Range("D43").Value = "0.282"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").Value = "2.804.66"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").Value = "38.15"
$ws.Range("E47").Value = "  -7.23%  "
$ws.Range("D48").Value = "128.47"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "24.59"
$ws.Range("E50").Value = "  -3.22%  "
$ws.Range("E51").Value = "  -1.44%  "
